# Auto-generated edit script: applies updated market-price values
# to the Typhon_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 225
$ws.Range("I33").Value = 238.21428
$ws.Range("K33").Value = 238.21428
$ws.Range("M33").Value = -9.214280000000002
$ws.Range("H116").Value = 22732754
$ws.Range("I116").Value = 250000000
$ws.Range("J116").Value = 6029.6
$ws.Range("K116").Value = 250000000
$ws.Range("L116").Value = 6029.6
$ws.Range("M116").Value = -249996558
$ws.Range("N116").Value = -12913.6
$ws.Range("H129").Value = 141770.53
$ws.Range("J129").Value = 164971.12
$ws.Range("L129").Value = 494913.36
$ws.Range("N129").Value = -504913.36
$ws.Range("H132").Value = 4158.75
$ws.Range("I132").Value = 4580.6665
$ws.Range("J132").Value = 1205.3334
$ws.Range("K132").Value = 13741.9995
$ws.Range("L132").Value = 3616.0002
$ws.Range("M132").Value = -11211.9995
$ws.Range("N132").Value = -8676.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6407.3945
$ws.Range("I32").Value = 4939.0483
$ws.Range("J32").Value = 16522.666
$ws.Range("K32").Value = 4939.0483
$ws.Range("L32").Value = 16522.666
$ws.Range("M32").Value = -4652.0483
$ws.Range("N32").Value = -17096.666
$ws.Range("H45").Value = 2126.7917
$ws.Range("I45").Value = 1905.6364
$ws.Range("J45").Value = 2313.923
$ws.Range("K45").Value = 1905.6364
$ws.Range("L45").Value = 2313.923
$ws.Range("M45").Value = -1528.6364
$ws.Range("N45").Value = -3067.923
$ws.Range("H74").Value = 19231934
$ws.Range("I74").Value = 23256320
$ws.Range("J74").Value = 4312.5557
$ws.Range("K74").Value = 23256320
$ws.Range("L74").Value = 4312.5557
$ws.Range("M74").Value = -23255446
$ws.Range("N74").Value = -6060.5557
$ws.Range("H77").Value = 19231934
$ws.Range("I77").Value = 23256320
$ws.Range("J77").Value = 4312.5557
$ws.Range("K77").Value = 116281600
$ws.Range("L77").Value = 21562.7785
$ws.Range("M77").Value = -116277232
$ws.Range("N77").Value = -30298.7785

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1434.4375
$ws.Range("J99").Value = 1329.1818
$ws.Range("L99").Value = 1329.1818
$ws.Range("N99").Value = -4325.1818
$ws.Range("H134").Value = 4053.64
$ws.Range("I134").Value = 4180.9165
$ws.Range("K134").Value = 12542.7495
$ws.Range("M134").Value = -10007.7495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1891.5
$ws.Range("I16").Value = 1891.5
$ws.Range("K16").Value = 1891.5
$ws.Range("M16").Value = -1604.5
$ws.Range("H31").Value = 4047.9714
$ws.Range("I31").Value = 2233.5293
$ws.Range("K31").Value = 2233.5293
$ws.Range("M31").Value = -1938.5293
$ws.Range("H34").Value = 4047.9714
$ws.Range("I34").Value = 2233.5293
$ws.Range("K34").Value = 2233.5293
$ws.Range("M34").Value = -2031.5293
$ws.Range("H99").Value = 2716.8572
$ws.Range("I99").Value = 2188.08
$ws.Range("K99").Value = 2188.08
$ws.Range("M99").Value = -690.0799999999999
$ws.Range("H113").Value = 1891.5
$ws.Range("I113").Value = 1891.5
$ws.Range("K113").Value = 1891.5
$ws.Range("M113").Value = 278.5
$ws.Range("H122").Value = 1767.5385
$ws.Range("I122").Value = 1859.7142
$ws.Range("J122").Value = 1660
$ws.Range("K122").Value = 5579.142599999999
$ws.Range("L122").Value = 4980
$ws.Range("M122").Value = -3129.142599999999
$ws.Range("N122").Value = -9880
$ws.Range("H126").Value = 2716.8572
$ws.Range("I126").Value = 2188.08
$ws.Range("K126").Value = 6564.24
$ws.Range("M126").Value = -4094.24
$ws.Range("H134").Value = 1341.0714
$ws.Range("I134").Value = 1214.6364
$ws.Range("K134").Value = 3643.9092
$ws.Range("M134").Value = -1108.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 9117.362999999999
$ws.Range("I2").Value = 10019.1
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 60114.60000000001
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -60001.60000000001
$ws.Range("N2").Value = -826
$ws.Range("H5").Value = 1317.9535
$ws.Range("J5").Value = 2043.2667
$ws.Range("L5").Value = 6129.800099999999
$ws.Range("N5").Value = -6353.800099999999
$ws.Range("H16").Value = 332.5
$ws.Range("I16").Value = 198.5
$ws.Range("J16").Value = 399.5
$ws.Range("K16").Value = 595.5
$ws.Range("L16").Value = 1198.5
$ws.Range("M16").Value = -422.5
$ws.Range("N16").Value = -1544.5
$ws.Range("H131").Value = 718.22
$ws.Range("J131").Value = 741.4194
$ws.Range("L131").Value = 2224.2582
$ws.Range("N131").Value = -12304.2582
$ws.Range("H135").Value = 1317.9535
$ws.Range("J135").Value = 2043.2667
$ws.Range("L135").Value = 18389.4003
$ws.Range("N135").Value = -23459.4003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 630.5
$ws.Range("I107").Value = 600
$ws.Range("J107").Value = 691.5
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 691.5
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -4531.5
$ws.Range("H113").Value = 9209.25
$ws.Range("I113").Value = 11145.667
$ws.Range("J113").Value = 3400
$ws.Range("K113").Value = 11145.667
$ws.Range("L113").Value = 3400
$ws.Range("M113").Value = -8975.666999999999
$ws.Range("N113").Value = -7740
$ws.Range("H132").Value = 23787.875
$ws.Range("I132").Value = 2381.6875
$ws.Range("J132").Value = 66600.25
$ws.Range("K132").Value = 7145.0625
$ws.Range("L132").Value = 199800.75
$ws.Range("M132").Value = -4615.0625
$ws.Range("N132").Value = -204860.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1703.4517
$ws.Range("I46").Value = 1907.2354
$ws.Range("J46").Value = 1456
$ws.Range("K46").Value = 1907.2354
$ws.Range("L46").Value = 1456
$ws.Range("M46").Value = -1719.2354
$ws.Range("N46").Value = -1832
$ws.Range("H55").Value = 930.2727
$ws.Range("I55").Value = 1082.3334
$ws.Range("K55").Value = 1082.3334
$ws.Range("M55").Value = -909.3334
$ws.Range("H68").Value = 2999.25
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2999.25
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2999.25
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4497.25
$ws.Range("H71").Value = 2999.25
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2999.25
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14996.25
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -22484.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 19500
$ws.Range("J28").Value = 19500
$ws.Range("L28").Value = 19500
$ws.Range("N28").Value = -20196
$ws.Range("H31").Value = 13000
$ws.Range("J31").Value = 13000
$ws.Range("L31").Value = 13000
$ws.Range("N31").Value = -13696
$ws.Range("H108").Value = 30312
$ws.Range("J108").Value = 30312
$ws.Range("L108").Value = 30312
$ws.Range("N108").Value = -37992
